# Update the 8.7.1 child-labour indicator sheet: relabel the "by sex" /
# "urban-rural" header rows and several English labels to their corrected
# (capitalised / reworded) forms, and fill in the previously-empty Kyrgyz
# label for "Functional difficulties in a child" in A35.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: "by sex" header row
$ws.Range("A6").Value = "Жынысы боюнча"
$ws.Range("B6").Value = "По полу"
$ws.Range("C6").Value = "By sex"

# Row 7: "Male"
$ws.Range("A7").Value = "Эркектер"
$ws.Range("B7").Value = "Мужчины"
$ws.Range("C7").Value = "Men"

# Row 8: "Female"
$ws.Range("A8").Value = "Аялдар"
$ws.Range("B8").Value = "Женщины"
$ws.Range("C8").Value = "Woman"

# Row 10: "urban"
$ws.Range("A10").Value = "Шаар"
$ws.Range("C10").Value = "Urban"

# Row 11: "rural"
$ws.Range("A11").Value = "Айыл"
$ws.Range("C11").Value = "Rural"

# Row 28: "does not attend"
$ws.Range("C28").Value = "Does not attend"

# Row 29: "educationof mother"
$ws.Range("C29").Value = "Educationof mother"

# Row 30: "preschool or not /primary"
$ws.Range("C30").Value = "Preschool or not /primary"

# Row 31: "basic general"
$ws.Range("C31").Value = "Basic general"

# Row 32: "average total"
$ws.Range("C32").Value = "Average total"

# Row 33: "vocational primary /secondary"
$ws.Range("C33").Value = "Vocational primary /secondary"

# Row 34: "higher"
$ws.Range("C34").Value = "Higher"

# Row 35: A35 was blank; fill in the Kyrgyz label matching B35/C35
$ws.Range("A35").Value = "Баланын функционалдык кыйнчылыктары"

# Row 38: "wealth quintile"
$ws.Range("C38").Value = "Wealth quintile"
